$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in row 17 (C17, D17, E17) which were previously empty
$ws.Range("C17").Value = "Avancement du code, fichier config, carte SD application, décodage NMEA, parsing fichier config."
$ws.Range("D17").Value = "-"
$ws.Range("E17").Value = 4

# Update selection / active cell to E18, removing the previous topLeftCell scroll position
$ws.Range("E18").Select()
